$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 1252, shifting all
# subsequent rows (and the final row 1351 -> 1352) down by one.
$ws.Rows.Item(1252).Insert()

# Populate the newly inserted row with its data.
$ws.Cells.Item(1252, 1).Value = 6
$ws.Cells.Item(1252, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1252, 3).Value = "Metropolitana"
$ws.Cells.Item(1252, 4).Value = 45106
$ws.Cells.Item(1252, 5).Value = 13
$ws.Cells.Item(1252, 6).Value = 100112021
$ws.Cells.Item(1252, 7).Value = "Ají"
$ws.Cells.Item(1252, 8).Value = "Americana (o)"
$ws.Cells.Item(1252, 9).Value = "Primera"
$ws.Cells.Item(1252, 10).Value = 500
$ws.Cells.Item(1252, 11).Value = 30000
$ws.Cells.Item(1252, 12).Value = 35000
$ws.Cells.Item(1252, 13).Value = 32700
$ws.Cells.Item(1252, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(1252, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(1252, 16).Value = 1308
$ws.Cells.Item(1252, 17).Value = 25
$ws.Cells.Item(1252, 18).Value = "Hortaliza"
